# Update the NATMI LR-pair result table (Sema3c -> Plxnd1) with the new
# TPM-derived values: existing sender/receiver combinations (ECs, FAPs,
# MuSCs) get refreshed statistics, and a new sending/receiving cluster
# "Resolving-Mac" is added, extending the sheet from 12 to 16 data rows
# (A2:T13 -> A2:T17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol,
# D=Target cluster, E..T = the various expression / specificity metrics.
$data = @(
  @("ECs", "Sema3c", "Plxnd1", "ECs", 3, 1, 0.604961, 1.814883, 0.0160517527720356, 0.0160517527720356, 3, 1, 139.2986196666667, 417.895859, 0.6137320738580456, 0.6137320738580456, 84.27023225216632, 758.4320902694969, 0.009851475517838042, 0.009851475517838043),
  @("ECs", "Sema3c", "Plxnd1", "FAPs", 3, 1, 0.604961, 1.814883, 0.0160517527720356, 0.0160517527720356, 3, 1, 22.17197066666667, 66.515912, 0.09768689432339951, 0.09768689432339951, 13.41317754647733, 120.718597918296, 0.001568045876747177, 0.001568045876747177),
  @("ECs", "Sema3c", "Plxnd1", "MuSCs", 3, 1, 0.604961, 1.814883, 0.0160517527720356, 0.0160517527720356, 3, 1, 35.78898466666666, 107.366954, 0.1576817331952585, 0.1576817331952585, 21.65093995293133, 194.858459576382, 0.002531068197916369, 0.00253106819791637),
  @("ECs", "Sema3c", "Plxnd1", "Resolving-Mac", 3, 1, 0.604961, 1.814883, 0.0160517527720356, 0.0160517527720356, 3, 1, 29.710182, 89.13054600000001, 0.1308992986232963, 0.1308992986232963, 17.973501412902, 161.761512716118, 0.002101163179534013, 0.002101163179534013),
  @("FAPs", "Sema3c", "Plxnd1", "ECs", 3, 1, 35.10542566666667, 105.316277, 0.9314709770686151, 0.9314709770686151, 3, 1, 139.2986196666667, 417.895859, 0.6137320738580456, 0.6137320738580456, 4890.137338177438, 44011.23604359694, 0.5716736144949013, 0.5716736144949013),
  @("FAPs", "Sema3c", "Plxnd1", "FAPs", 3, 1, 35.10542566666667, 105.316277, 0.9314709770686151, 0.9314709770686151, 3, 1, 22.17197066666667, 66.515912, 0.09768689432339951, 0.09768689432339951, 778.3564681221804, 7005.208213099624, 0.0909925069022155, 0.0909925069022155),
  @("FAPs", "Sema3c", "Plxnd1", "MuSCs", 3, 1, 35.10542566666667, 105.316277, 0.9314709770686151, 0.9314709770686151, 3, 1, 35.78898466666666, 107.366954, 0.1576817331952585, 0.1576817331952585, 1256.38754090114, 11307.48786811026, 0.1468759580852601, 0.1468759580852602),
  @("FAPs", "Sema3c", "Plxnd1", "Resolving-Mac", 3, 1, 35.10542566666667, 105.316277, 0.9314709770686151, 0.9314709770686151, 3, 1, 29.710182, 89.13054600000001, 0.1308992986232963, 0.1308992986232963, 1042.988585744138, 9386.897271697242, 0.1219288975862383, 0.1219288975862383),
  @("MuSCs", "Sema3c", "Plxnd1", "ECs", 3, 1, 1.958375, 5.875125, 0.05196260806057782, 0.05196260806057783, 3, 1, 139.2986196666667, 417.895859, 0.6137320738580456, 0.6137320738580456, 272.7989342897083, 2455.190408607375, 0.03189111920809123, 0.03189111920809123),
  @("MuSCs", "Sema3c", "Plxnd1", "FAPs", 3, 1, 1.958375, 5.875125, 0.05196260806057782, 0.05196260806057783, 3, 1, 22.17197066666667, 66.515912, 0.09768689432339951, 0.09768689432339951, 43.42103305433334, 390.789297489, 0.005076065802381893, 0.005076065802381894),
  @("MuSCs", "Sema3c", "Plxnd1", "MuSCs", 3, 1, 1.958375, 5.875125, 0.05196260806057782, 0.05196260806057783, 3, 1, 35.78898466666666, 107.366954, 0.1576817331952585, 0.1576817331952585, 70.08825284658333, 630.7942756192499, 0.008193554100337822, 0.008193554100337824),
  @("MuSCs", "Sema3c", "Plxnd1", "Resolving-Mac", 3, 1, 1.958375, 5.875125, 0.05196260806057782, 0.05196260806057783, 3, 1, 29.710182, 89.13054600000001, 0.1308992986232963, 0.1308992986232963, 58.18367767425001, 523.6530990682501, 0.006801868949766881, 0.006801868949766882),
  @("Resolving-Mac", "Sema3c", "Plxnd1", "ECs", 1, 0.3333333333333333, 0.01939666666666667, 0.05819, 0.0005146620987715195, 0.0005146620987715195, 3, 1, 139.2986196666667, 417.895859, 0.6137320738580456, 0.6137320738580456, 2.701928892801111, 24.31736003521, 0.000315864637215179, 0.000315864637215179),
  @("Resolving-Mac", "Sema3c", "Plxnd1", "FAPs", 1, 0.3333333333333333, 0.01939666666666667, 0.05819, 0.0005146620987715195, 0.0005146620987715195, 3, 1, 22.17197066666667, 66.515912, 0.09768689432339951, 0.09768689432339951, 0.4300623243644444, 3.87056091928, [double]"5.027574205495243e-05", [double]"5.027574205495243e-05"),
  @("Resolving-Mac", "Sema3c", "Plxnd1", "MuSCs", 1, 0.3333333333333333, 0.01939666666666667, 0.05819, 0.0005146620987715195, 0.0005146620987715195, 3, 1, 35.78898466666666, 107.366954, 0.1576817331952585, 0.1576817331952585, 0.6941870059177777, 6.247683053259999, [double]"8.115281174420253e-05", [double]"8.115281174420255e-05"),
  @("Resolving-Mac", "Sema3c", "Plxnd1", "Resolving-Mac", 1, 0.3333333333333333, 0.01939666666666667, 0.05819, 0.0005146620987715195, 0.0005146620987715195, 3, 1, 29.710182, 89.13054600000001, 0.1308992986232963, 0.1308992986232963, 0.5762784968600001, 5.18650647174, [double]"6.736890775718556e-05", [double]"6.736890775718556e-05")
)

# Write the table starting at row 2 (row 1 holds the headers and is
# left untouched), columns A(1) through T(20).
$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value2 = $val
        $c++
    }
    $r++
}

Write-Host "Rows written: $($data.Length), dimension now A1:T$($r - 1)"
